# Update "想去人数" (want-to-go count) figures on the "展览" and "全部类型" sheets.
$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (exhibition) ---
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F2").Value = 369
$wsExhibit.Range("F4").Value = 279
$wsExhibit.Range("F5").Value = 4129
$wsExhibit.Range("F6").Value = 40
$wsExhibit.Range("F7").Value = 456

# --- Sheet "全部类型" (all types) ---
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F2").Value = 369
$wsAll.Range("F4").Value = 279
$wsAll.Range("F5").Value = 4129
$wsAll.Range("F8").Value = 40
$wsAll.Range("F9").Value = 456
